$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/17_kotako2.wav"
$ws.Range("B2").Value = "pngimages/17_cracker.png"

$ws.Range("A3").Value = "trainingaudio/25_tapapi1.wav"
$ws.Range("B3").Value = "pngimages/25_apple.png"

$ws.Range("A4").Value = "trainingaudio/12_pokika3.wav"
$ws.Range("B4").Value = "pngimages/12_pie.png"

$ws.Range("A5").Value = "trainingaudio/07_pitapi2.wav"
$ws.Range("B5").Value = "pngimages/07_suitcase.png"

$ws.Range("A6").Value = "trainingaudio/11_tokiko1.wav"
$ws.Range("B6").Value = "pngimages/11_compass.png"

$ws.Range("A7").Value = "trainingaudio/27_pakapa1.wav"
$ws.Range("B7").Value = "pngimages/27_kiwi.png"
